$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Work Plan v2.3 - split Task 5.5.2 "Text to Speech Conversion" effort so
# that the three existing language-specific sub-tasks (English / isiZulu /
# Kinyarwanda, rows 47-49) each drop from 0.5 to 0.25 person-months, and a
# new fourth sub-task "Integrated Text to Speech Conversion" is inserted
# right after them (new row 50) with the freed-up 0.25 person-months.
# Every row from the old row 50 onward shifts down by one row.
# ---------------------------------------------------------------------------

# 1) Insert a new blank row at 50 - this pushes rows 50.. down to 51.. and
#    lets Excel auto-adjust every formula/range reference that spans the
#    insertion point (SUM ranges, shared-formula ranges, totals, etc).
$ws.Rows("50:50").Insert()

# 2) Copy the formatting (styles) of the row that is now 51 (old row 50)
#    onto the newly inserted row 50, so the new row looks identical to its
#    neighbours instead of using Excel's generic blank-row formatting.
$ws.Range("A51:G51").Copy()
$ws.Range("A50:G50").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Halve the effort fraction for the three existing language sub-tasks;
#    column D (Task Duration = B*C) recalculates automatically.
$ws.Range("C47").Value = 0.25
$ws.Range("C48").Value = 0.25
$ws.Range("C49").Value = 0.25

# 4) Fill in the new sub-task row (row 50) with the "Integrated" task.
$bom = [char]0xFEFF
$ws.Range("A50").Value = $bom + "Task 5.5.2.2 Integrated Text to Speech Conversion"
$ws.Range("B50").Value = 6
$ws.Range("C50").Value = 0.25
$ws.Range("D50").Formula = "=SUM(B50*C50)"
$ws.Range("G50").Formula = "=SUM(`$D50)"

# 5) Match the author's final selection in the saved file.
$ws.Range("C51").Select()
